# Auto-generated edit script applying the cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.17"
$ws.Range("D3").Value = "'35.78"
$ws.Range("E3").Value = "'0.43%"
$ws.Range("D4").Value = "'5.127"
$ws.Range("E4").Value = "'1.10%"
$ws.Range("D5").Value = "'0.08113"
$ws.Range("E5").Value = "'3.19%"
$ws.Range("D6").Value = "'2.129"
$ws.Range("E6").Value = "'0.04%"
$ws.Range("E7").Value = "'1.37%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9306"
$ws.Range("E8").Value = "'1.23%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1017"
$ws.Range("E9").Value = "'4.24%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1870"
$ws.Range("E10").Value = "'0.75%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09224"
$ws.Range("E11").Value = "'6.79%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03614"
$ws.Range("E12").Value = "'1.84%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09891"
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001438"
$ws.Range("E14").Value = "'0.68%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005794"
$ws.Range("E15").Value = "'2.89%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.468"
$ws.Range("E16").Value = "'0.28%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.150"
$ws.Range("E17").Value = "'1.18%"
$ws.Range("D18").Value = "'2.816"
$ws.Range("E18").Value = "'7.58%"
$ws.Range("D19").Value = "'0.3369"
$ws.Range("E19").Value = "'-1.69%"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("E20").Value = "'1.25%"
$ws.Range("D21").Value = "'5.131"
$ws.Range("E21").Value = "'-1.78%"
$ws.Range("D23").Value = "'0.04583"
$ws.Range("E23").Value = "'0.71%"
$ws.Range("E24").Value = "'1.16%"
$ws.Range("E25").Value = "'-6.89%"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("E26").Value = "'-21.81%"
$ws.Range("D27").Value = "'0.0004520"
$ws.Range("E27").Value = "'-4.84%"
$ws.Range("D39").Value = "'0.01965"
$ws.Range("E39").Value = "'6.33%"
$ws.Range("D40").Value = "'0.04879"
$ws.Range("E40").Value = "'3.30%"
$ws.Range("D41").Value = "'0.007826"
$ws.Range("E41").Value = "'4.38%"
$ws.Range("D42").Value = "'0.1389"
$ws.Range("E42").Value = "'-0.52%"
$ws.Range("D43").Value = "'0.007872"
$ws.Range("E43").Value = "'1.60%"
$ws.Range("D44").Value = "'0.002109"
$ws.Range("E44").Value = "'-5.86%"
$ws.Range("D45").Value = "'0.01164"
$ws.Range("E45").Value = "'3.12%"
$ws.Range("D46").Value = "'0.00006527"
$ws.Range("E46").Value = "'2.99%"
$ws.Range("E47").Value = "'0.43%"
$ws.Range("D48").Value = "'39.12"
$ws.Range("E48").Value = "'-16.73%"
$ws.Range("D49").Value = "'0.001909"
$ws.Range("E49").Value = "'-4.58%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.43%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.43%"
